# Generate Report for Handback
#
# The file 8a1b7149-7071-4066-8814-3de76465f0f2.md has now been handed
# back (in sync with en-US) for both the zh-cn and de-de locales.
# This re-sorts each sheet so that the handed-back file moves to row 2
# (previously row 3), pushing the still-in-flight 2ed6ece3-... file to
# row 3, and it populates the "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" columns for the handed-back file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Common URL fragments re-used for the various hyperlinks
# ---------------------------------------------------------------
$md8a1bUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/49187da92123e31175b6420c8931f5b0a05037ae/e2e/8a1b7149-7071-4066-8814-3de76465f0f2.md"
$md2ed6Url   = "https://github.com/OpenLocalizationTest/oltest/blob/8725b1402a4143afdeeddc3440056647a8b80dff/e2e/2ed6ece3-4e41-4e90-bc16-b993d4daa480.md"
$xlf8a1bZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/903879448f07b8dd4665685736ff73ddabfdedb8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8a1b7149-7071-4066-8814-3de76465f0f2.7051c9a158f5f74d62c702eb7d967f02488d84bc.zh-cn.xlf"
$xlf2ed6ZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f7f5a860bbbf32e6c6b9c9bfe38790afb1c5760/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2ed6ece3-4e41-4e90-bc16-b993d4daa480.313c817be3bf3bd53310360001f135646b44bb13.zh-cn.xlf"
$xlf8a1bDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/06985e41841096c6c6f09e1120b01fe2d8120f89/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8a1b7149-7071-4066-8814-3de76465f0f2.7051c9a158f5f74d62c702eb7d967f02488d84bc.de-de.xlf"
$xlf2ed6DeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a38c16d0201bd014453348630b7eb63a0f9b1f71/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2ed6ece3-4e41-4e90-bc16-b993d4daa480.313c817be3bf3bd53310360001f135646b44bb13.de-de.xlf"

$md8a1bName   = "8a1b7149-7071-4066-8814-3de76465f0f2.md"
$md2ed6Name   = "2ed6ece3-4e41-4e90-bc16-b993d4daa480.md"
$xlf8a1bZhName = "8a1b7149-7071-4066-8814-3de76465f0f2.7051c9a158f5f74d62c702eb7d967f02488d84bc.zh-cn.xlf"
$xlf2ed6ZhName = "2ed6ece3-4e41-4e90-bc16-b993d4daa480.313c817be3bf3bd53310360001f135646b44bb13.zh-cn.xlf"
$xlf8a1bDeName = "8a1b7149-7071-4066-8814-3de76465f0f2.7051c9a158f5f74d62c702eb7d967f02488d84bc.de-de.xlf"
$xlf2ed6DeName = "2ed6ece3-4e41-4e90-bc16-b993d4daa480.313c817be3bf3bd53310360001f135646b44bb13.de-de.xlf"

$handedBack  = "Handed back: in sync with en-US"
$readyStatus = "Ready for handoff"

# =================================================================
# Sheet "Overview"
# =================================================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $md8a1bName
$ov.Range("B2").Value = $handedBack
$ov.Range("C2").Value = $handedBack
$ov.Range("D2").Value = "2016-27-14 03:27:35"

$ov.Range("A3").Value = $md2ed6Name
$ov.Range("B3").Value = $readyStatus
$ov.Range("C3").Value = $readyStatus
$ov.Range("D3").Value = "2016-27-14 03:27:21"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), $md8a1bUrl, [Type]::Missing, [Type]::Missing, $md8a1bName) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), $md2ed6Url, [Type]::Missing, [Type]::Missing, $md2ed6Name) | Out-Null

# =================================================================
# Sheet "zh-cn"
# =================================================================
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $md8a1bName
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $handedBack
$zh.Range("D2").Value = $xlf8a1bZhName
$zh.Range("E2").Value = "2016-03-14 03:27:33"
$zh.Range("F2").Value = $md8a1bName
$zh.Range("G2").Value = $xlf8a1bZhName
$zh.Range("H2").Value = "2016-03-14 03:27:49"
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = $md2ed6Name
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = $readyStatus
$zh.Range("D3").Value = $xlf2ed6ZhName
$zh.Range("E3").Value = "2016-03-14 03:27:18"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("I3").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $md8a1bUrl, [Type]::Missing, [Type]::Missing, $md8a1bName) | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), $md8a1bUrl, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), $xlf8a1bZhUrl, [Type]::Missing, [Type]::Missing, $xlf8a1bZhName) | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), $md8a1bUrl, [Type]::Missing, [Type]::Missing, $md8a1bName) | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), $xlf8a1bZhUrl, [Type]::Missing, [Type]::Missing, $xlf8a1bZhName) | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), $md2ed6Url, [Type]::Missing, [Type]::Missing, $md2ed6Name) | Out-Null
$zh.Hyperlinks.Add($zh.Range("B3"), $md2ed6Url, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), $xlf2ed6ZhUrl, [Type]::Missing, [Type]::Missing, $xlf2ed6ZhName) | Out-Null

# =================================================================
# Sheet "de-de"
# =================================================================
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $md8a1bName
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $handedBack
$de.Range("D2").Value = $xlf8a1bDeName
$de.Range("E2").Value = "2016-03-14 03:27:35"
$de.Range("F2").Value = $md8a1bName
$de.Range("G2").Value = $xlf8a1bDeName
$de.Range("H2").Value = "2016-03-14 03:27:54"
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = $md2ed6Name
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = $readyStatus
$de.Range("D3").Value = $xlf2ed6DeName
$de.Range("E3").Value = "2016-03-14 03:27:21"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("I3").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $md8a1bUrl, [Type]::Missing, [Type]::Missing, $md8a1bName) | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), $md8a1bUrl, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), $xlf8a1bDeUrl, [Type]::Missing, [Type]::Missing, $xlf8a1bDeName) | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), $md8a1bUrl, [Type]::Missing, [Type]::Missing, $md8a1bName) | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), $xlf8a1bDeUrl, [Type]::Missing, [Type]::Missing, $xlf8a1bDeName) | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), $md2ed6Url, [Type]::Missing, [Type]::Missing, $md2ed6Name) | Out-Null
$de.Hyperlinks.Add($de.Range("B3"), $md2ed6Url, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), $xlf2ed6DeUrl, [Type]::Missing, [Type]::Missing, $xlf2ed6DeName) | Out-Null

Write-Host "Handback report generated."
